# Apply the edits described by the diff to the active workbook/worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update / insert the data rows ---
# Original data (before edit):
#   Row2: Cololabis saira | Pacific Saury | Teleost Fish
#   Row3: Unassigned      | Unassigned    | Unassigned
#
# New data (after edit):
#   Row2: Fundulus heteroclitus or majalis | Mummichog or striped killifish | Teleost Fish
#   Row3: Cololabis saira                  | Pacific suary                  | Teleost Fish
#   Row4: Unassigned                       | Unassigned                     | Unassigned
#   Row5: Mareca americana                 | American wigeon                | Bird
#   Row6: Myrophis vafer                   | Pacific worm eel               | Teleost Fish

$ws.Range("A2").Value = "Fundulus heteroclitus or majalis"
$ws.Range("B2").Value = "Mummichog or striped killifish"
$ws.Range("C2").Value = "Teleost Fish"

$ws.Range("A3").Value = "Cololabis saira"
$ws.Range("B3").Value = "Pacific suary"
$ws.Range("C3").Value = "Teleost Fish"

$ws.Range("A4").Value = "Unassigned"
$ws.Range("B4").Value = "Unassigned"
$ws.Range("C4").Value = "Unassigned"

$ws.Range("A5").Value = "Mareca americana"
$ws.Range("B5").Value = "American wigeon"
$ws.Range("C5").Value = "Bird"

$ws.Range("A6").Value = "Myrophis vafer"
$ws.Range("B6").Value = "Pacific worm eel"
$ws.Range("C6").Value = "Teleost Fish"

# --- Widen the three columns (23.77734375 -> 34.5546875 stored width) ---
# The ColumnWidth setter in this engine rounds to the nearest 1/6th of a
# character, so 33.6667 is the input that lands closest to the recorded
# target width of 34.5546875 (rounds to 34.5, the nearest reachable step).
$ws.Columns("A:C").ColumnWidth = 33.666666666667

# --- Update the active selection to B4 (was A4) ---
$ws.Range("B4").Select()
